# Insert a new data row at row 131, pushing the existing rows 131..184
# down to 132..185 (mirrors the OOXML diff: dimension A1:R184 -> A1:R185).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(131).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(131, 1).Value  = 3
$ws.Cells.Item(131, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(131, 3).Value  = "Coquimbo"
$ws.Cells.Item(131, 4).Value  = 44795
$ws.Cells.Item(131, 5).Value  = 5
$ws.Cells.Item(131, 6).Value  = 100112026
$ws.Cells.Item(131, 7).Value  = "Haba"
$ws.Cells.Item(131, 8).Value  = "Sin especificar"
$ws.Cells.Item(131, 9).Value  = "Primera"
$ws.Cells.Item(131, 10).Value = 65
$ws.Cells.Item(131, 11).Value = 14000
$ws.Cells.Item(131, 12).Value = 14000
$ws.Cells.Item(131, 13).Value = 14000
$ws.Cells.Item(131, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(131, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(131, 16).Value = 560
$ws.Cells.Item(131, 17).Value = 25
$ws.Cells.Item(131, 18).Value = "Hortaliza"
